$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 381 (shifts existing rows 381..436 down to 382..437)
$ws.Rows.Item(381).Insert()

# Populate the newly inserted row 381 with the new weekly data point
$ws.Cells.Item(381, 1).Value = 9
$ws.Cells.Item(381, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(381, 3).Value = "Metropolitana"
$ws.Cells.Item(381, 4).Value = 45127
$ws.Cells.Item(381, 5).Value = 13
$ws.Cells.Item(381, 6).Value = 100112043
$ws.Cells.Item(381, 7).Value = "Pepino ensalada"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 70
$ws.Cells.Item(381, 11).Value = 13000
$ws.Cells.Item(381, 12).Value = 14000
$ws.Cells.Item(381, 13).Value = 13500
$ws.Cells.Item(381, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(381, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(381, 16).Value = 270
$ws.Cells.Item(381, 17).Value = 50
$ws.Cells.Item(381, 18).Value = "Hortaliza"
